$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "828÷8=103, 4" "939÷6=156, 3"
Replace-Text "979÷3=326, 1" "974÷6=162, 2"
Replace-Text "777÷3=259, 0" "140÷8=17, 4"
Replace-Text "271÷2=135, 1" "769÷5=153, 4"
Replace-Text "909÷5=181, 4" "266÷6=44, 2"
Replace-Text "755÷9=83, 8" "681÷6=113, 3"
Replace-Text "444÷2=222, 0" "933÷9=103, 6"
Replace-Text "533÷6=88, 5" "574÷6=95, 4"
Replace-Text "625÷7=89, 2" "102÷3=34, 0"
Replace-Text "879÷7=125, 4" "787÷6=131, 1"
Replace-Text "422÷8=52, 6" "376÷4=94, 0"
Replace-Text "331÷7=47, 2" "454÷4=113, 2"
Replace-Text "543÷9=60, 3" "250÷6=41, 4"
Replace-Text "380÷9=42, 2" "517÷2=258, 1"
Replace-Text "450÷6=75, 0" "908÷4=227, 0"
Replace-Text "708÷2=354, 0" "586÷2=293, 0"
Replace-Text "806÷7=115, 1" "627÷2=313, 1"
Replace-Text "404÷3=134, 2" "310÷4=77, 2"
Replace-Text "147÷9=16, 3" "766÷2=383, 0"
Replace-Text "307÷5=61, 2" "714÷9=79, 3"
Replace-Text "439÷6=73, 1" "432÷6=72, 0"
Replace-Text "531÷2=265, 1" "825÷5=165, 0"
Replace-Text "279÷6=46, 3" "646÷2=323, 0"
Replace-Text "303÷8=37, 7" "436÷3=145, 1"
Replace-Text "427÷8=53, 3" "436÷5=87, 1"

Write-Output "Replacements applied"
